{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that reads \"LOB1019: F\u00edsica II (Requisito)\"; the\n// footer block to remove (a blank paragraph, the \"Ver no Jupiter...\" line,\n// and the \"\u00a9 2020 ...\" copyright line) immediately follows it.\nconst items = paragraphs.items;\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"F\u00edsica II (Requisito)\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOB1019: F\u00edsica II (Requisito)' paragraph\");\n}\n\n// Delete the next three paragraphs: blank line, \"Ver no Jupiter...\" line,\n// and the \"\u00a9 2020 ...\" copyright line.\nfor (let i = 0; i < 3; i++) {\n  items[anchorIndex + 1 + i].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph \"LOB1019: F\u00edsica II (Requisito)\". The footer block to\n# remove (a blank paragraph, the \"Ver no Jupiter...\" line, and the\n# \"\u00a9 2020 ...\" copyright line) immediately follows it.\n$anchor = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*F\u00edsica II (Requisito)*\") {\n        $anchor = $i\n        break\n    }\n}\n\nif ($anchor -eq -1) {\n    throw \"Could not find the 'LOB1019: F\u00edsica II (Requisito)' paragraph\"\n}\n\n# Delete the next three paragraphs: blank line, \"Ver no Jupiter...\" line,\n# and the \"\u00a9 2020 ...\" copyright line. Re-fetching Item($anchor + 1) each\n# time accounts for the collection shrinking after each delete.\nfor ($j = 0; $j -lt 3; $j++) {\n    $d.Paragraphs.Item($anchor + 1).Range.Delete()\n}\n"}
